$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '64.470.41'
$ws.Range('E2').Value2 = '  +1.46%  '
$ws.Range('D3').Value2 = '3.188.18'
$ws.Range('E3').Value2 = '  +2.86%  '
$ws.Range('E4').Value2 = '  +0.19%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '593.84'
$c.Style = "Normal"
$ws.Range('E5').Value2 = '  +1.77%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '149.03'
$c.Style = "Normal"
$ws.Range('E7').Value2 = '  +0.13%  '
$ws.Range('D8').Value2 = '3.177.33'
$ws.Range('E8').Value2 = '  +2.73%  '
$ws.Range('E9').Value2 = '  +1.25%  '
$ws.Range('E10').Value2 = '  +1.65%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '6.01'
$c.Style = "Normal"
$ws.Range('E11').Value2 = '  +6.70%  '
$ws.Range('E12').Value2 = '  +1.43%  '
$ws.Range('E13').Value2 = '  +1.21%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '37.90'
$c.Style = "Normal"
$ws.Range('E14').Value2 = '  +2.20%  '
$ws.Range('D15').Value2 = '3.722.21'
$ws.Range('E15').Value2 = '  +3.09%  '
$ws.Range('E16').Value2 = '  +0.14%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '7.39'
$c.Style = "Normal"
$ws.Range('E17').Value2 = '  +4.60%  '
$ws.Range('D18').Value2 = '3.188.70'
$ws.Range('E18').Value2 = '  +2.92%  '
$ws.Range('D19').Value2 = '64.261.17'
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '474.25'
$c.Style = "Normal"
$ws.Range('E20').Value2 = '  +3.01%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '14.63'
$c.Style = "Normal"
$ws.Range('E21').Value2 = '  +2.99%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.741'
$c.Style = "Normal"
$ws.Range('E22').Value2 = '  +2.63%  '
$ws.Range('E23').Value2 = '  +3.79%  '
$ws.Range('E24').Value2 = '  +9.28%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '13.34'
$c.Style = "Normal"
$ws.Range('E25').Value2 = '  +3.30%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '81.96'
$c.Style = "Normal"
$ws.Range('E26').Value2 = '  +1.12%  '
$ws.Range('B27').Value2 = 'Dai'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E27').Value2 = '  +0.08%  '
$ws.Range('B28').Value2 = 'RenderToken'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '10.00'
$c.Style = "Normal"
$ws.Range('E28').Value2 = '  +8.81%  '
$ws.Range('B29').Value2 = 'PancakeSwap'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.74'
$c.Style = "Normal"
$ws.Range('E29').Value2 = '  +2.86%  '
$ws.Range('B30').Value2 = 'ImmutableX'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.26'
$c.Style = "Normal"
$ws.Range('E30').Value2 = '  +3.00%  '
$ws.Range('E31').Value2 = '  +0.19%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '7.27'
$c.Style = "Normal"
$ws.Range('E32').Value2 = '  +4.45%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.121'
$c.Style = "Normal"
$ws.Range('E33').Value2 = '  +9.38%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '28.51'
$c.Style = "Normal"
$ws.Range('E34').Value2 = '  +6.98%  '
$ws.Range('D35').Value2 = '0.0₃0866'
$ws.Range('E35').Value2 = '  +2.05%  '
$ws.Range('E36').Value2 = '  +3.79%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '6.26'
$c.Style = "Normal"
$ws.Range('E37').Value2 = '  +4.48%  '
$ws.Range('E38').Value2 = '  +1.36%  '
$ws.Range('E39').Value2 = '  +0.57%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '469.26'
$c.Style = "Normal"
$ws.Range('E40').Value2 = '  +7.63%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '9.46'
$c.Style = "Normal"
$ws.Range('E41').Value2 = '  +9.02%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '51.56'
$c.Style = "Normal"
$ws.Range('E42').Value2 = '  +2.59%  '
$ws.Range('E43').Value2 = '  +7.96%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.0379'
$c.Style = "Normal"
$ws.Range('E44').Value2 = '  +3.17%  '
$ws.Range('D45').Value2 = '2.940.55'
$ws.Range('E45').Value2 = '  +1.98%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '39.48'
$c.Style = "Normal"
$ws.Range('E46').Value2 = '  +8.30%  '
$ws.Range('E47').Value2 = '  +1.34%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '132.99'
$c.Style = "Normal"
$ws.Range('E48').Value2 = '  +6.56%  '
$ws.Range('E49').Value2 = '  +0.05%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.28'
$c.Style = "Normal"
$ws.Range('E50').Value2 = '  +6.13%  '
$ws.Range('E51').Value2 = '  +1.54%  '
